$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.278.25'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '2.612.62'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.87'
$ws.Range('E5').Value = '  +2.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.25'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.597'
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('E10').Value = '  -0.59%  '
$ws.Range('E11').Value = '  +2.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.156'
$ws.Range('E12').Value = '  +1.20%  '
$ws.Range('D13').Value = '3.073.68'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.22'
$ws.Range('E14').Value = '  +7.38%  '
$ws.Range('D15').Value = '60.273.82'
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').Value = '2.615.66'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.45'
$ws.Range('E18').Value = '  +2.21%  '
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '346.71'
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('E21').Value = '  -2.43%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  +1.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.68'
$ws.Range('E24').Value = '  -0.69%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.03'
$ws.Range('E27').Value = '  +4.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.94'
$ws.Range('E28').Value = '  +5.86%  '
$ws.Range('D29').Value = '0.0₃0799'
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.51'
$ws.Range('E30').Value = '  +2.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '168.90'
$ws.Range('E31').Value = '  +4.20%  '
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  +5.92%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.30'
$ws.Range('E35').Value = '  +1.44%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.30'
$ws.Range('E36').Value = '  +7.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.63'
$ws.Range('E37').Value = '  +1.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '319.87'
$ws.Range('E38').Value = '  +7.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '38.50'
$ws.Range('E39').Value = '  +1.99%  '
$ws.Range('E40').Value = '  +4.18%  '
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '135.67'
$ws.Range('E42').Value = '  -3.09%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.99'
$ws.Range('E44').Value = '  +2.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0989'
$ws.Range('E45').Value = '  +0.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.610'
$ws.Range('E46').Value = '  +0.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.05'
$ws.Range('E47').Value = '  +3.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0552'
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('E49').Value = '  +1.79%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('E51').Value = '  +0.50%  '
